$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new worksheet column at K (11th column) so everything from
# "Doc 1 Final Diagnosis" onward shifts one column to the right, opening
# up a slot for the new "AI Confidence" column right after
# "AI Recommendation".
$ws.Columns.Item(11).Insert()

# Grow the table range so it covers the newly inserted column.
$lo.Resize($ws.Range("D2:O12"))

# Header + data for the new "AI Confidence" column.
$ws.Range("K2").Value = "AI Confidence"
$ws.Range("K3").Value = 0.94
$ws.Range("K4").Value = 0.6
$ws.Range("K5").Value = 0.8
$ws.Range("K6").Value = "…"
$ws.Range("K7").Value = "-"
$ws.Range("K8").Value = "-"
$ws.Range("K9").Value = "…"
$ws.Range("K10").Value = 0.7
$ws.Range("K11").Value = 0.86
$ws.Range("K12").Value = "…"
